$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reset row 3 (previously filled-in Yumi robot bug report) back to the blank
# placeholder template text, matching the pattern used in row 2 -- except
# the status cell (H3) keeps reading "Error".
$ws.Range("A3").Value = "[Where or which part the error or bug found]"
$ws.Range("B3").Value = "[Describe the bug.]"
$ws.Range("C3").Value = "[Describe the expected result.]"
$ws.Range("D3").Value = "[Enter the actual result. ]"
$ws.Range("E3").Value = "[Person who wrote this test case.]"
$ws.Range("F3").Value = "[Date when this test case was last updated.]"
$ws.Range("H3").Value = "Error"

# F3 held a real date value (styled with a date number format); now that it
# holds placeholder text again it should look like the rest of the row
# (same formatting as F2), so copy F2's formatting onto it.
$ws.Range("F2").Copy()
$ws.Range("F3").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 3 shrinks back down from the tall, data-filled height to a compact
# template-row height.
$ws.Rows.Item(3).RowHeight = 64.8

# Selection moves from the single filled-in cell to the header+template
# block.
$ws.Range("A1:H3").Select() | Out-Null
